# Update the cryptocurrency price/volume table on Sheet1 to the latest
# scraped values (GitHub Actions refresh commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.913.04'
$ws.Range('E2').Value = '  -0.48%  '
$ws.Range('D3').Value = '1.667.72'
$ws.Range('E3').Value = '  +0.87%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.23'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.520'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.08%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.0623'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.37%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.249'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.33'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.81%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0892'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.86%  '
$ws.Range('D12').Value = '1.901.75'
$ws.Range('E12').Value = '  +0.82%  '
$ws.Range('D13').Value = '1.661.24'
$ws.Range('E14').Value = '  +0.11%  '
$ws.Range('E15').Value = '  +1.46%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.80'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.70%  '
$ws.Range('D17').Value = '26.902.30'
$ws.Range('E17').Value = '  -0.54%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '235.12'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.95'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.45%  '
$ws.Range('E20').Value = '  +0.27%  '
$ws.Range('E21').Value = '  +0.00%  '
$ws.Range('E22').Value = '  -0.50%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.17'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.87%  '
$ws.Range('E24').Value = '  -2.93%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.69'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.12'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('E27').Value = '  -1.05%  '
$ws.Range('E28').Value = '  +0.15%  '
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0494'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.71%  '
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('E32').Value = '  +1.09%  '
$ws.Range('D33').Value = '1.445.78'
$ws.Range('E33').Value = '  -5.18%  '
$ws.Range('E34').Value = '  +1.47%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.64'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.91%  '
$ws.Range('E36').Value = '  -0.19%  '
$ws.Range('E37').Value = '  +0.56%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.904'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.75%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0171'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.81%  '
$ws.Range('E40').Value = '  -3.78%  '
$ws.Range('E41').Value = '  +0.11%  '
$ws.Range('B42').Value = 'WEMIXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +9.51%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.30'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '65.98'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.41%  '
$ws.Range('D45').Value = '1.808.39'
$ws.Range('E45').Value = '  +0.86%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.782'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.13%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.75'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.34%  '
$ws.Range('E48').Value = '  +1.04%  '
$ws.Range('E49').Value = '  -1.18%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.102'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.16%  '
$ws.Range('E51').Value = '  -0.10%  '
